$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns stay as text (avoid Excel auto-numeric conversion)

# Row 45: Aave -> WEMIXToken
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.917"
$ws.Range("E45").Value = "  -10.21%  "

# Row 46: WEMIXToken -> Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.74"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.401.08"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.614.73"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.99"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0609"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.65"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.613.06"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.65"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.427.54"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.96"
$ws.Range("E18").Value = "  +6.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.16"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.03"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.447.66"
$ws.Range("E33").Value = "  +8.54%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  -5.13%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.753.43"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.92"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("E48").Value = "  +4.55%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -1.26%  "
